$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UOVEY")

# Apply updated financial figures (Income Statement, Balance Sheet, Cash Flow Statement)
$ws.Range("D8").Value = 6692600
$ws.Range("E8").Value = 6113000
$ws.Range("F8").Value = 5770600
$ws.Range("G8").Value = 5301100
$ws.Range("H8").Value = 4798800
$ws.Range("I8").Value = 4572800
$ws.Range("J8").Value = 4159000

$ws.Range("D17").Value = 3152800
$ws.Range("E17").Value = 2870900
$ws.Range("F17").Value = 2633500
$ws.Range("G17").Value = 2408800
$ws.Range("H17").Value = 2077300
$ws.Range("I17").Value = 2035800
$ws.Range("J17").Value = 1832900

$ws.Range("D18").Value = 3539700
$ws.Range("E18").Value = 3242100
$ws.Range("F18").Value = 3137100
$ws.Range("G18").Value = 2892200
$ws.Range("H18").Value = 2721500
$ws.Range("I18").Value = 2537000
$ws.Range("J18").Value = 2326200

$ws.Range("D20").Value = -437900
$ws.Range("E20").Value = -457200
$ws.Range("F20").Value = -284500
$ws.Range("G20").Value = -71900
$ws.Range("H20").Value = -78800
$ws.Range("I20").Value = -66000
$ws.Range("J20").Value = -255700

$ws.Range("D21").Value = 3292500
$ws.Range("E21").Value = 2948500
$ws.Range("F21").Value = 2986500
$ws.Range("G21").Value = 2940800
$ws.Range("H21").Value = 2738700
$ws.Range("I21").Value = 2560600
$ws.Range("J21").Value = "NA"

$ws.Range("D23").Value = 3101900
$ws.Range("E23").Value = 2784900
$ws.Range("F23").Value = 2852600
$ws.Range("G23").Value = 2820300
$ws.Range("H23").Value = 2642800
$ws.Range("I23").Value = 2471000
$ws.Range("J23").Value = 2070400

$ws.Range("D24").Value = 590000
$ws.Range("E24").Value = 493300
$ws.Range("F24").Value = 478300
$ws.Range("G24").Value = 413400
$ws.Range("H24").Value = 412200
$ws.Range("I24").Value = 391300
$ws.Range("J24").Value = 344500

$ws.Range("D26").Value = 2511900
$ws.Range("E26").Value = 2291600
$ws.Range("F26").Value = 2374300
$ws.Range("G26").Value = 2406900
$ws.Range("H26").Value = 2230500
$ws.Range("I26").Value = 2079700
$ws.Range("J26").Value = 1725900

$ws.Range("D27").Value = 2424400
$ws.Range("E27").Value = 2214900
$ws.Range("F27").Value = 2288400
$ws.Range("G27").Value = 2320400
$ws.Range("H27").Value = 2141800
$ws.Range("I27").Value = 1991000
$ws.Range("J27").Value = 1639900

$ws.Range("D32").Value = 437900
$ws.Range("E32").Value = 457200
$ws.Range("F32").Value = 284500
$ws.Range("G32").Value = 71900
$ws.Range("H32").Value = 78800
$ws.Range("I32").Value = 66000
$ws.Range("J32").Value = 255700

$ws.Range("D33").Value = 2424400
$ws.Range("E33").Value = 2214900
$ws.Range("F33").Value = 2288400
$ws.Range("G33").Value = 2320400
$ws.Range("H33").Value = 2141800
$ws.Range("I33").Value = 1991000
$ws.Range("J33").Value = 1639900

$ws.Range("D35").Value = 2424400
$ws.Range("E35").Value = 2214900
$ws.Range("F35").Value = 2288400
$ws.Range("G35").Value = 2320400
$ws.Range("H35").Value = 2141800
$ws.Range("I35").Value = 1991000
$ws.Range("J35").Value = 1639900

$ws.Range("D41").Value = 77739500
$ws.Range("E41").Value = 47452200
$ws.Range("F41").Value = 44942600
$ws.Range("G41").Value = 47024500
$ws.Range("H41").Value = 42981800
$ws.Range("I41").Value = 36164800
$ws.Range("J41").Value = 19750800

$ws.Range("D42").Value = 9614300
$ws.Range("E42").Value = 10888500
$ws.Range("F42").Value = 8729500
$ws.Range("G42").Value = 8979300
$ws.Range("H42").Value = 11372300
$ws.Range("I42").Value = 8397700
$ws.Range("J42").Value = 200100

$ws.Range("D47").Value = 880400
$ws.Range("E47").Value = 817700
$ws.Range("F47").Value = 815500
$ws.Range("G47").Value = 877000
$ws.Range("H47").Value = 734800
$ws.Range("I47").Value = 812600
$ws.Range("J47").Value = 805400

$ws.Range("D48").Value = 4511200
$ws.Range("E48").Value = 2204500
$ws.Range("F48").Value = 2099300
$ws.Range("G48").Value = 1761100
$ws.Range("H48").Value = 1691000
$ws.Range("I48").Value = 3317800
$ws.Range("J48").Value = 1604600

$ws.Range("D49").Value = 6108100
$ws.Range("E49").Value = 3060800
$ws.Range("F49").Value = 3055400
$ws.Range("G49").Value = 3059500
$ws.Range("H49").Value = 3055400
$ws.Range("I49").Value = 6146800
$ws.Range("J49").Value = 3093900

$ws.Range("D52").Value = 142500
$ws.Range("E52").Value = 185100
$ws.Range("F52").Value = 167400
$ws.Range("G52").Value = 170800
$ws.Range("H52").Value = 212100
$ws.Range("I52").Value = 218800
$ws.Range("J52").Value = 246200

$ws.Range("D54").Value = 264408000
$ws.Range("E54").Value = 250719000
$ws.Range("F54").Value = 233011000
$ws.Range("G54").Value = 226172000
$ws.Range("H54").Value = 209576000
$ws.Range("I54").Value = 186476000
$ws.Range("J54").Value = 174721000

$ws.Range("D57").Value = 2629000
$ws.Range("E57").Value = 3083500
$ws.Range("F57").Value = 2676600
$ws.Range("G57").Value = 1453700
$ws.Range("H57").Value = 1401600
$ws.Range("I57").Value = 393300
$ws.Range("J57").Value = 292600

$ws.Range("D59").Value = 1098700
$ws.Range("E59").Value = 850000
$ws.Range("F59").Value = 844400
$ws.Range("G59").Value = 759600
$ws.Range("H59").Value = 833000
$ws.Range("I59").Value = 886600
$ws.Range("J59").Value = 805000

$ws.Range("D61").Value = 18565000
$ws.Range("E61").Value = 19276500
$ws.Range("F61").Value = 14959600
$ws.Range("G61").Value = 15449900
$ws.Range("H61").Value = 13995900
$ws.Range("I61").Value = 9438100
$ws.Range("J61").Value = 8690400

$ws.Range("D62").Value = 131200
$ws.Range("E62").Value = 171000
$ws.Range("F62").Value = 141300
$ws.Range("G62").Value = 118300
$ws.Range("H62").Value = 63700
$ws.Range("I62").Value = 16000
$ws.Range("J62").Value = 25200

$ws.Range("D66").Value = 237237000
$ws.Range("E66").Value = 226480000
$ws.Range("F66").Value = 210324000
$ws.Range("G66").Value = 204369000
$ws.Range("H66").Value = 190119000
$ws.Range("I66").Value = 167983000
$ws.Range("J66").Value = 157786000

$ws.Range("F70").Value = 613100
$ws.Range("G70").Value = 613100
$ws.Range("H70").Value = 613100

$ws.Range("D72").Value = 28370400
$ws.Range("E72").Value = 19891600
$ws.Range("F72").Value = 18498400
$ws.Range("G72").Value = 17519300
$ws.Range("H72").Value = 16056900
$ws.Range("I72").Value = 22555800
$ws.Range("J72").Value = 13697400

$ws.Range("D76").Value = 27171300
$ws.Range("E76").Value = 24239000
$ws.Range("F76").Value = 22074000
$ws.Range("G76").Value = 21189900
$ws.Range("H76").Value = 18844000
$ws.Range("I76").Value = 18492700
$ws.Range("J76").Value = 16934800

$ws.Range("D81").Value = 2424400
$ws.Range("E81").Value = 2214900
$ws.Range("F81").Value = 2288400
$ws.Range("G81").Value = 2320400
$ws.Range("H81").Value = 2141800
$ws.Range("I81").Value = 1991000
$ws.Range("J81").Value = 1639900

$ws.Range("D83").Value = 190400
$ws.Range("E83").Value = 163500
$ws.Range("F83").Value = 133800
$ws.Range("G83").Value = 120500
$ws.Range("H83").Value = 95900
$ws.Range("I83").Value = 89600
$ws.Range("J83").Value = "NA"

$ws.Range("D89").Value = 2882300
$ws.Range("E89").Value = -9816700
$ws.Range("F89").Value = 487500
$ws.Range("G89").Value = 5243200
$ws.Range("H89").Value = -8623400
$ws.Range("I89").Value = 4674800
$ws.Range("J89").Value = -7048700

$ws.Range("D91").Value = -257600
$ws.Range("E91").Value = -282300
$ws.Range("F91").Value = -510200
$ws.Range("G91").Value = -190700
$ws.Range("H91").Value = -163200
$ws.Range("I91").Value = -162600
$ws.Range("J91").Value = -138000

$ws.Range("D94").Value = -251400
$ws.Range("E94").Value = -281600
$ws.Range("F94").Value = -356500
$ws.Range("G94").Value = 44700
$ws.Range("H94").Value = -52400
$ws.Range("I94").Value = -50000
$ws.Range("J94").Value = "NA"

$ws.Range("D96").Value = -490300
$ws.Range("E96").Value = -687200
$ws.Range("F96").Value = -1093200
$ws.Range("G96").Value = -521800
$ws.Range("H96").Value = -889000
$ws.Range("I96").Value = -771900
$ws.Range("J96").Value = -474000

$ws.Range("D100").Value = -259600
$ws.Range("E100").Value = 3482000
$ws.Range("F100").Value = -1726100
$ws.Range("G100").Value = 883100
$ws.Range("H100").Value = 3666800
$ws.Range("I100").Value = -36500
$ws.Range("J100").Value = "NA"

$ws.Range("D101").Value = -473000
$ws.Range("E101").Value = 107600
$ws.Range("F101").Value = -230800
$ws.Range("G101").Value = 67200
$ws.Range("H101").Value = -188700
$ws.Range("I101").Value = -240000
$ws.Range("J101").Value = "NA"

$ws.Range("D102").Value = 1898300
$ws.Range("E102").Value = -6508600
$ws.Range("F102").Value = -1826000
$ws.Range("G102").Value = 6238200
$ws.Range("H102").Value = -5197600
$ws.Range("I102").Value = 4348400
$ws.Range("J102").Value = -3500100

